# Aggiorna test_turni.xlsx con nuovi parametri
$wb = $excel.ActiveWorkbook

# --- Sheet "Pianificazione": reshuffle some shift assignments ---
$ws1 = $wb.Worksheets.Item("Pianificazione")

# Row 12 (09/11/2025): remove Luigi Bianchi from Turno Pomeriggio (D12)
$ws1.Range("D12").Value = ""

# Row 19 (16/11/2025): remove Luigi Bianchi from Turno Mattina (C19)
$ws1.Range("C19").Value = ""

# Row 26 (23/11/2025): remove Luigi Bianchi from Turno Pomeriggio (D26)
$ws1.Range("D26").Value = ""

# Row 28 (25/11/2025): add Luigi Bianchi to Turno Pomeriggio (D28)
$ws1.Range("D28").Value = "Luigi Bianchi"

# Row 29 (26/11/2025): add Luigi Bianchi to Turno Mattina (C29)
$ws1.Range("C29").Value = "Luigi Bianchi"

# Row 30 (27/11/2025): add Luigi Bianchi to Turno Pomeriggio (D30)
$ws1.Range("D30").Value = "Luigi Bianchi"

# Row 31 (28/11/2025): add Luigi Bianchi to Turno Mattina (C31) and Mario Rossi to Turno Pomeriggio (D31)
$ws1.Range("C31").Value = "Luigi Bianchi"
$ws1.Range("D31").Value = "Mario Rossi"

# Row 32 (29/11/2025): add Mario Rossi to Turno Mattina (C32) and Luigi Bianchi to Turno Pomeriggio (D32)
$ws1.Range("C32").Value = "Mario Rossi"
$ws1.Range("D32").Value = "Luigi Bianchi"

# Row 33 (30/11/2025): add Mario Rossi to Turno Pomeriggio (D33)
$ws1.Range("D33").Value = "Mario Rossi"

# --- Sheet "Statistiche": update computed totals ---
$ws2 = $wb.Worksheets.Item("Statistiche")

$ws2.Range("B5").Value = 180
$ws2.Range("E5").Value = 30

$ws2.Range("B6").Value = 156
$ws2.Range("E6").Value = 26

$ws2.Range("B15").Value = 5
$ws2.Range("B16").Value = 1

# --- Sheet "Dettagli Addetti": rename headers and update max-hours values ---
$ws3 = $wb.Worksheets.Item("Dettagli Addetti")

$ws3.Range("B3").Value = "Ore Contratto (min)"
$ws3.Range("C3").Value = "Ore Max (sett)"

$ws3.Range("C4").Value = 45
$ws3.Range("C5").Value = 40
